# Generate Report for Archive
#
# 1. Status text "Ready for handoff" -> "In Translation" on all three sheets
#    (Overview!E2:F3, zh-cn!C2:C3, de-de!C2:C3).
# 2. Narrow the "Status" columns (Overview E:F, zh-cn C, de-de C) from
#    ~17.22 chars to ~13.41 chars.

$wb = $excel.ActiveWorkbook

$overview = $wb.Worksheets.Item("Overview")
$zhcn     = $wb.Worksheets.Item("zh-cn")
$dede     = $wb.Worksheets.Item("de-de")

# --- Update the "Status" text for every row on every sheet ---
$overview.Cells.Replace("Ready for handoff", "In Translation") | Out-Null
$zhcn.Cells.Replace("Ready for handoff", "In Translation") | Out-Null
$dede.Cells.Replace("Ready for handoff", "In Translation") | Out-Null

# --- Narrow the Status columns ---
$overview.Columns.Item(5).ColumnWidth = 12.5   # column E (zh-cn status)
$overview.Columns.Item(6).ColumnWidth = 12.5   # column F (de-de status)
$zhcn.Columns.Item(3).ColumnWidth = 12.5        # column C (Status)
$dede.Columns.Item(3).ColumnWidth = 12.5        # column C (Status)
